$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order swapped fixture rows (home/away pairs whose scrape order changed) ---
# Row 3
$ws.Range("F3").Value = "Konyaspor"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "Istanbulspor AS"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1.94
$ws.Range("K3").Value = "06/08/2023 21:12"
$ws.Range("L3").Value = 1.83
$ws.Range("M3").Value = "12/08/2023 18:06"
$ws.Range("N3").Value = 3.83
$ws.Range("O3").Value = "06/08/2023 21:12"
$ws.Range("P3").Value = 3.9
$ws.Range("Q3").Value = "12/08/2023 18:06"
$ws.Range("R3").Value = 3.87
$ws.Range("S3").Value = "06/08/2023 21:12"
$ws.Range("T3").Value = 4.41
$ws.Range("U3").Value = "12/08/2023 18:06"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/turkey/super-lig/konyaspor-istanbulspor-as/06o0xHmB/"

# Row 4
$ws.Range("F4").Value = "Kasimpasa"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = "Ankaragucu"
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 2.33
$ws.Range("K4").Value = "07/08/2023 01:12"
$ws.Range("L4").Value = 2.73
$ws.Range("M4").Value = "12/08/2023 18:14"
$ws.Range("N4").Value = 3.67
$ws.Range("O4").Value = "07/08/2023 01:12"
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = "12/08/2023 18:11"
$ws.Range("R4").Value = 3.03
$ws.Range("S4").Value = "07/08/2023 01:12"
$ws.Range("T4").Value = 2.67
$ws.Range("U4").Value = "12/08/2023 18:13"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-ankaragucu/URAVMylO/"

# Row 18
$ws.Range("F18").Value = "Besiktas"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "Pendikspor"
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1.29
$ws.Range("K18").Value = "15/08/2023 13:42"
$ws.Range("L18").Value = 1.28
$ws.Range("M18").Value = "20/08/2023 20:12"
$ws.Range("N18").Value = 6.05
$ws.Range("O18").Value = "15/08/2023 13:42"
$ws.Range("P18").Value = 6.31
$ws.Range("Q18").Value = "20/08/2023 20:12"
$ws.Range("R18").Value = 8.81
$ws.Range("S18").Value = "15/08/2023 13:42"
$ws.Range("T18").Value = 10.07
$ws.Range("U18").Value = "20/08/2023 20:12"
$ws.Range("V18").Value = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-pendikspor/6RGYjG24/"

# Row 19
$ws.Range("F19").Value = "Gaziantep"
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = "Sivasspor"
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 2.45
$ws.Range("K19").Value = "15/08/2023 13:42"
$ws.Range("L19").Value = 2.84
$ws.Range("M19").Value = "20/08/2023 20:43"
$ws.Range("N19").Value = 3.59
$ws.Range("O19").Value = "15/08/2023 13:42"
$ws.Range("P19").Value = 3.42
$ws.Range("Q19").Value = "20/08/2023 20:43"
$ws.Range("R19").Value = 2.9
$ws.Range("S19").Value = "15/08/2023 13:42"
$ws.Range("T19").Value = 2.61
$ws.Range("U19").Value = "20/08/2023 20:43"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-sivasspor/KAEMgIIo/"

# Row 31
$ws.Range("F31").Value = "Basaksehir"
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = "Konyaspor"
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 2.06
$ws.Range("K31").Value = "26/08/2023 18:43"
$ws.Range("L31").Value = 2.27
$ws.Range("M31").Value = "02/09/2023 18:14"
$ws.Range("N31").Value = 3.66
$ws.Range("O31").Value = "26/08/2023 18:43"
$ws.Range("P31").Value = 3.48
$ws.Range("Q31").Value = "02/09/2023 18:14"
$ws.Range("R31").Value = 3.65
$ws.Range("S31").Value = "26/08/2023 18:43"
$ws.Range("T31").Value = 3.34
$ws.Range("U31").Value = "02/09/2023 18:14"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-konyaspor/f7gp56NE/"

# Row 32
$ws.Range("F32").Value = "Pendikspor"
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = "Alanyaspor"
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = 2.47
$ws.Range("K32").Value = "28/08/2023 04:12"
$ws.Range("L32").Value = 2.78
$ws.Range("M32").Value = "02/09/2023 18:13"
$ws.Range("N32").Value = 3.74
$ws.Range("O32").Value = "28/08/2023 04:12"
$ws.Range("P32").Value = 3.53
$ws.Range("Q32").Value = "02/09/2023 18:14"
$ws.Range("R32").Value = 2.73
$ws.Range("S32").Value = "28/08/2023 04:12"
$ws.Range("T32").Value = 2.6
$ws.Range("U32").Value = "02/09/2023 18:13"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/turkey/super-lig/pendikspor-alanyaspor/6oDkqP0e/"

# Row 91
$ws.Range("F91").Value = "Hatayspor"
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = "Kayserispor"
$ws.Range("I91").Value = 2
$ws.Range("J91").Value = 2.26
$ws.Range("K91").Value = "22/10/2023 20:15"
$ws.Range("L91").Value = 2.49
$ws.Range("M91").Value = "27/10/2023 18:56"
$ws.Range("N91").Value = 3.55
$ws.Range("O91").Value = "22/10/2023 20:15"
$ws.Range("P91").Value = 3.42
$ws.Range("Q91").Value = "27/10/2023 18:56"
$ws.Range("R91").Value = 3.26
$ws.Range("S91").Value = "22/10/2023 20:15"
$ws.Range("T91").Value = 3.01
$ws.Range("U91").Value = "27/10/2023 18:56"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-kayserispor/069EREiq/"

# Row 92
$ws.Range("F92").Value = "Kasimpasa"
$ws.Range("G92").Value = 3
$ws.Range("H92").Value = "Istanbulspor AS"
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = 1.86
$ws.Range("K92").Value = "22/10/2023 15:12"
$ws.Range("L92").Value = 1.96
$ws.Range("M92").Value = "27/10/2023 18:58"
$ws.Range("N92").Value = 3.96
$ws.Range("O92").Value = "22/10/2023 15:12"
$ws.Range("P92").Value = 3.71
$ws.Range("Q92").Value = "27/10/2023 18:59"
$ws.Range("R92").Value = 4.09
$ws.Range("S92").Value = "22/10/2023 15:12"
$ws.Range("T92").Value = 3.99
$ws.Range("U92").Value = "27/10/2023 18:58"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-istanbulspor-as/YJ7MPhMe/"

# Row 97
$ws.Range("F97").Value = "Antalyaspor"
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = "Basaksehir"
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2.03
$ws.Range("K97").Value = "22/10/2023 20:15"
$ws.Range("L97").Value = 2.11
$ws.Range("M97").Value = "29/10/2023 16:54"
$ws.Range("N97").Value = 3.51
$ws.Range("O97").Value = "22/10/2023 20:15"
$ws.Range("P97").Value = 3.35
$ws.Range("Q97").Value = "29/10/2023 16:54"
$ws.Range("R97").Value = 3.8
$ws.Range("S97").Value = "22/10/2023 20:15"
$ws.Range("T97").Value = 3.89
$ws.Range("U97").Value = "29/10/2023 16:54"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-basaksehir/tz8GozqF/"

# Row 98
$ws.Range("F98").Value = "Ankaragucu"
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = "Samsunspor"
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2.19
$ws.Range("K98").Value = "23/10/2023 05:42"
$ws.Range("L98").Value = 2.7
$ws.Range("M98").Value = "29/10/2023 16:59"
$ws.Range("N98").Value = 3.59
$ws.Range("O98").Value = "23/10/2023 05:42"
$ws.Range("P98").Value = 3.32
$ws.Range("Q98").Value = "29/10/2023 16:54"
$ws.Range("R98").Value = 3.29
$ws.Range("S98").Value = "23/10/2023 05:42"
$ws.Range("T98").Value = 2.82
$ws.Range("U98").Value = "29/10/2023 16:59"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-samsunspor/2kdPqEDR/"

# Row 99
$ws.Range("F99").Value = "Pendikspor"
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = "Fenerbahce"
$ws.Range("I99").Value = 5
$ws.Range("J99").Value = 7.87
$ws.Range("K99").Value = "23/10/2023 05:42"
$ws.Range("L99").Value = 10.56
$ws.Range("M99").Value = "29/10/2023 16:59"
$ws.Range("N99").Value = 5.68
$ws.Range("O99").Value = "23/10/2023 05:42"
$ws.Range("P99").Value = 6.22
$ws.Range("Q99").Value = "29/10/2023 16:59"
$ws.Range("R99").Value = 1.35
$ws.Range("S99").Value = "23/10/2023 05:42"
$ws.Range("T99").Value = 1.28
$ws.Range("U99").Value = "29/10/2023 16:59"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/turkey/super-lig/pendikspor-fenerbahce/vc8IQY6k/"

# Row 112
$ws.Range("F112").Value = "Trabzonspor"
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = "Konyaspor"
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = 1.76
$ws.Range("K112").Value = "05/11/2023 14:12"
$ws.Range("L112").Value = 1.62
$ws.Range("M112").Value = "10/11/2023 17:53"
$ws.Range("N112").Value = 3.95
$ws.Range("O112").Value = "05/11/2023 14:12"
$ws.Range("P112").Value = 4.04
$ws.Range("Q112").Value = "10/11/2023 17:46"
$ws.Range("R112").Value = 4.66
$ws.Range("S112").Value = "05/11/2023 14:12"
$ws.Range("T112").Value = 6.11
$ws.Range("U112").Value = "10/11/2023 17:46"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/turkey/super-lig/trabzonspor-konyaspor/QoWcf4Rn/"

# Row 113
$ws.Range("F113").Value = "Ankaragucu"
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = "Antalyaspor"
$ws.Range("I113").Value = 4
$ws.Range("J113").Value = 2.4
$ws.Range("K113").Value = "05/11/2023 17:12"
$ws.Range("L113").Value = 2.47
$ws.Range("M113").Value = "10/11/2023 17:55"
$ws.Range("N113").Value = 3.45
$ws.Range("O113").Value = "05/11/2023 17:12"
$ws.Range("P113").Value = 3.33
$ws.Range("Q113").Value = "10/11/2023 17:55"
$ws.Range("R113").Value = 3.07
$ws.Range("S113").Value = "05/11/2023 17:12"
$ws.Range("T113").Value = 3.1
$ws.Range("U113").Value = "10/11/2023 17:55"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-antalyaspor/MDP6hrda/"

# Row 119
$ws.Range("F119").Value = "Besiktas"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = "Basaksehir"
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1.55
$ws.Range("K119").Value = "05/11/2023 17:12"
$ws.Range("L119").Value = 1.74
$ws.Range("M119").Value = "12/11/2023 13:56"
$ws.Range("N119").Value = 4.44
$ws.Range("O119").Value = "05/11/2023 17:12"
$ws.Range("P119").Value = 3.84
$ws.Range("Q119").Value = "12/11/2023 13:59"
$ws.Range("R119").Value = 5.79
$ws.Range("S119").Value = "05/11/2023 17:12"
$ws.Range("T119").Value = 5.16
$ws.Range("U119").Value = "12/11/2023 13:59"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-basaksehir/vkz2gOtg/"

# Row 120
$ws.Range("F120").Value = "Rizespor"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Istanbulspor AS"
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 1.8
$ws.Range("K120").Value = "06/11/2023 18:12"
$ws.Range("L120").Value = 1.69
$ws.Range("M120").Value = "12/11/2023 13:53"
$ws.Range("N120").Value = 3.94
$ws.Range("O120").Value = "06/11/2023 18:12"
$ws.Range("P120").Value = 4.01
$ws.Range("Q120").Value = "12/11/2023 13:53"
$ws.Range("R120").Value = 4.32
$ws.Range("S120").Value = "06/11/2023 18:12"
$ws.Range("T120").Value = 5.3
$ws.Range("U120").Value = "12/11/2023 13:52"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/turkey/super-lig/rizespor-istanbulspor-as/04OAi2B5/"

# Row 127
$ws.Range("F127").Value = "Samsunspor"
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = "Besiktas"
$ws.Range("I127").Value = 2
$ws.Range("J127").Value = 3.24
$ws.Range("K127").Value = "12/11/2023 15:42"
$ws.Range("L127").Value = 3.31
$ws.Range("M127").Value = "26/11/2023 13:57"
$ws.Range("N127").Value = 3.56
$ws.Range("O127").Value = "12/11/2023 15:42"
$ws.Range("P127").Value = 3.5
$ws.Range("Q127").Value = "26/11/2023 13:57"
$ws.Range("R127").Value = 2.23
$ws.Range("S127").Value = "12/11/2023 15:42"
$ws.Range("T127").Value = 2.27
$ws.Range("U127").Value = "26/11/2023 13:57"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-besiktas/UcVQmKAU/"

# Row 128
$ws.Range("F128").Value = "Antalyaspor"
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = "Rizespor"
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 1.67
$ws.Range("K128").Value = "12/11/2023 15:42"
$ws.Range("L128").Value = 1.68
$ws.Range("M128").Value = "26/11/2023 13:57"
$ws.Range("N128").Value = 4.06
$ws.Range("O128").Value = "12/11/2023 15:42"
$ws.Range("P128").Value = 3.88
$ws.Range("Q128").Value = "26/11/2023 13:55"
$ws.Range("R128").Value = 5.07
$ws.Range("S128").Value = "12/11/2023 15:42"
$ws.Range("T128").Value = 5.55
$ws.Range("U128").Value = "26/11/2023 13:55"
$ws.Range("V128").Value = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-rizespor/2o0u7sBh/"

# --- Append two new fixture rows scraped in this run ---
# Row 130
$ws.Range("A130").Value = 129
$ws.Range("B130").Value = "turkey"
$ws.Range("C130").Value = "super-lig"
$ws.Range("D130").Value = "2023-2024"
$ws.Range("E130").Value = 45257.75
$ws.Range("F130").Value = "Basaksehir"
$ws.Range("G130").Value = 4
$ws.Range("H130").Value = "Pendikspor"
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = 1.74
$ws.Range("K130").Value = "12/11/2023 15:42"
$ws.Range("L130").Value = 2.16
$ws.Range("M130").Value = "27/11/2023 17:59"
$ws.Range("N130").Value = 3.95
$ws.Range("O130").Value = "12/11/2023 15:42"
$ws.Range("P130").Value = 3.49
$ws.Range("Q130").Value = "27/11/2023 17:59"
$ws.Range("R130").Value = 4.9
$ws.Range("S130").Value = "12/11/2023 15:42"
$ws.Range("T130").Value = 3.58
$ws.Range("U130").Value = "27/11/2023 17:59"
$ws.Range("V130").Value = "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-pendikspor/pA5z8Nen/"

# Row 131
$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "turkey"
$ws.Range("C131").Value = "super-lig"
$ws.Range("D131").Value = "2023-2024"
$ws.Range("E131").Value = 45257.75
$ws.Range("F131").Value = "Sivasspor"
$ws.Range("G131").Value = 3
$ws.Range("H131").Value = "Trabzonspor"
$ws.Range("I131").Value = 3
$ws.Range("J131").Value = 2.91
$ws.Range("K131").Value = "12/11/2023 15:42"
$ws.Range("L131").Value = 3.3
$ws.Range("M131").Value = "27/11/2023 17:59"
$ws.Range("N131").Value = 3.38
$ws.Range("O131").Value = "12/11/2023 15:42"
$ws.Range("P131").Value = 3.38
$ws.Range("Q131").Value = "27/11/2023 17:59"
$ws.Range("R131").Value = 2.55
$ws.Range("S131").Value = "12/11/2023 15:42"
$ws.Range("T131").Value = 2.33
$ws.Range("U131").Value = "27/11/2023 17:57"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/turkey/super-lig/sivasspor-trabzonspor/zabm5Lu5/"

# --- Match formatting used by the other data rows (bold/bordered index column, date-time column) ---
$idx = $ws.Range("A130")
$idx.Font.Bold = $true
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160
$idx.Borders.LineStyle = 1
$dt = $ws.Range("E130")
$dt.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$idx = $ws.Range("A131")
$idx.Font.Bold = $true
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160
$idx.Borders.LineStyle = 1
$dt = $ws.Range("E131")
$dt.NumberFormat = "YYYY-MM-DD HH:MM:SS"
